# update scripts wuth new tpm
# Re-run of the NATMI ligand-receptor (Timp3 -> Agtr2) edge table with
# updated TPM-normalized expression values, plus a new "Resolving-Mac"
# sending cluster (2 additional rows, for Target cluster FAPs/MuSCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2: ECs -> Timp3 | Agtr2 -> FAPs ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Timp3"
$ws.Range("C2").Value = "Agtr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 68.88366733333334
$ws.Range("H2").Value = 206.651002
$ws.Range("I2").Value = 0.5393713802555014
$ws.Range("J2").Value = 0.5393713802555014
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.002414333333333
$ws.Range("N2").Value = 6.007243
$ws.Range("O2").Value = 0.7140239834365498
$ws.Range("P2").Value = 0.7140239834365498
$ws.Range("Q2").Value = 137.9336428008318
$ws.Range("R2").Value = 1241.402785207486
$ws.Range("S2").Value = 0.3851241014817031
$ws.Range("T2").Value = 0.3851241014817031

# --- row 3: ECs -> Timp3 | Agtr2 -> MuSCs ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Timp3"
$ws.Range("C3").Value = "Agtr2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 68.88366733333334
$ws.Range("H3").Value = 206.651002
$ws.Range("I3").Value = 0.5393713802555014
$ws.Range("J3").Value = 0.5393713802555014
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8019933333333333
$ws.Range("N3").Value = 2.40598
$ws.Range("O3").Value = 0.2859760165634502
$ws.Range("P3").Value = 0.2859760165634502
$ws.Range("Q3").Value = 55.24424197688445
$ws.Range("R3").Value = 497.19817779196
$ws.Range("S3").Value = 0.1542472787737982
$ws.Range("T3").Value = 0.1542472787737982

# --- row 4: FAPs -> Timp3 | Agtr2 -> FAPs ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Timp3"
$ws.Range("C4").Value = "Agtr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 27.97197233333334
$ws.Range("H4").Value = 83.91591700000001
$ws.Range("I4").Value = 0.2190255239009008
$ws.Range("J4").Value = 0.2190255239009008
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.002414333333333
$ws.Range("N4").Value = 6.007243
$ws.Range("O4").Value = 0.7140239834365498
$ws.Range("P4").Value = 0.7140239834365498
$ws.Range("Q4").Value = 56.01147833187012
$ws.Range("R4").Value = 504.1033049868311
$ws.Range("S4").Value = 0.1563894770499984
$ws.Range("T4").Value = 0.1563894770499984

# --- row 5: FAPs -> Timp3 | Agtr2 -> MuSCs ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Timp3"
$ws.Range("C5").Value = "Agtr2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 27.97197233333334
$ws.Range("H5").Value = 83.91591700000001
$ws.Range("I5").Value = 0.2190255239009008
$ws.Range("J5").Value = 0.2190255239009008
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8019933333333333
$ws.Range("N5").Value = 2.40598
$ws.Range("O5").Value = 0.2859760165634502
$ws.Range("P5").Value = 0.2859760165634502
$ws.Range("Q5").Value = 22.43333533151778
$ws.Range("R5").Value = 201.90001798366
$ws.Range("S5").Value = 0.06263604685090235
$ws.Range("T5").Value = 0.06263604685090235

# --- row 6: MuSCs -> Timp3 | Agtr2 -> FAPs ---
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Timp3"
$ws.Range("C6").Value = "Agtr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 30.64425733333333
$ws.Range("H6").Value = 91.932772
$ws.Range("I6").Value = 0.2399499912628263
$ws.Range("J6").Value = 0.2399499912628263
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.002414333333333
$ws.Range("N6").Value = 6.007243
$ws.Range("O6").Value = 0.7140239834365498
$ws.Range("P6").Value = 0.7140239834365498
$ws.Range("Q6").Value = 61.36250011862177
$ws.Range("R6").Value = 552.262501067596
$ws.Range("S6").Value = 0.1713300485870486
$ws.Range("T6").Value = 0.1713300485870486

# --- row 7: MuSCs -> Timp3 | Agtr2 -> MuSCs ---
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Timp3"
$ws.Range("C7").Value = "Agtr2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 30.64425733333333
$ws.Range("H7").Value = 91.932772
$ws.Range("I7").Value = 0.2399499912628263
$ws.Range("J7").Value = 0.2399499912628263
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8019933333333333
$ws.Range("N7").Value = 2.40598
$ws.Range("O7").Value = 0.2859760165634502
$ws.Range("P7").Value = 0.2859760165634502
$ws.Range("Q7").Value = 24.57649008628444
$ws.Range("R7").Value = 221.18841077656
$ws.Range("S7").Value = 0.06861994267577774
$ws.Range("T7").Value = 0.06861994267577776

# --- row 8 (new): Resolving-Mac -> Timp3 | Agtr2 -> FAPs ---
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Timp3"
$ws.Range("C8").Value = "Agtr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2111196666666667
$ws.Range("H8").Value = 0.633359
$ws.Range("I8").Value = 0.001653104580771614
$ws.Range("J8").Value = 0.001653104580771615
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.002414333333333
$ws.Range("N8").Value = 6.007243
$ws.Range("O8").Value = 0.7140239834365498
$ws.Range("P8").Value = 0.7140239834365498
$ws.Range("Q8").Value = 0.4227490465818889
$ws.Range("R8").Value = 3.804741419237
$ws.Range("S8").Value = 0.001180356317799756
$ws.Range("T8").Value = 0.001180356317799756

# --- row 9 (new): Resolving-Mac -> Timp3 | Agtr2 -> MuSCs ---
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Timp3"
$ws.Range("C9").Value = "Agtr2"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2111196666666667
$ws.Range("H9").Value = 0.633359
$ws.Range("I9").Value = 0.001653104580771614
$ws.Range("J9").Value = 0.001653104580771615
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8019933333333333
$ws.Range("N9").Value = 2.40598
$ws.Range("O9").Value = 0.2859760165634502
$ws.Range("P9").Value = 0.2859760165634502
$ws.Range("Q9").Value = 0.1693165652022222
$ws.Range("R9").Value = 1.52384908682
$ws.Range("S9").Value = 0.0004727482629718586
$ws.Range("T9").Value = 0.0004727482629718586
